# Updates cryptos list values (price + volume%) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.742.08"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "1.567.51"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "25.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "1.790.01"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "1.564.41"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "28.725.57"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "0.0₃0681"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.105"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0461"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.63%  "
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "1.398.42"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("E35").Value = "  -3.65%  "
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.21%  "
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0162"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0457"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").Value = "1.703.05"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.871"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0513"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.21%  "
